$d = $word.ActiveDocument

# 1. Set the row height (trHeight) on the second row of the last table (row with
#    "Added work" / documentation cell) to 692 twips (= 34.6 points).
$t = $d.Tables.Item($d.Tables.Count)
$row = $t.Rows.Item(2)
$row.Height = 34.6

# 2. Update the documentation bullet text.
$d.Content.Find.Execute(
    "Documentatie " + [char]8211 + " Terminat Arhitectura aplicatiei (Onion + 3-Tier)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Documentatie -Terminat capitolul Proiectare si implementare", 2) | Out-Null

# 3. Remove the "De facut decriptarea obligatorie la medicament" run while keeping
#    the bookmark that follows it in the same paragraph.
$d.Content.Find.Execute(
    "De facut decriptarea obligatorie la medicament",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null
